$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each block below rewrites the cells of one data row (D, K, L, M, N, O, P,
# Q, R, S, T) to match the updated weekly Damasco price report. The data
# turned out to be a reshuffle of the 16 existing rows (rows 2-17) plus a
# couple of incidental value tweaks, so every row below is fully
# re-stated with its final values.

# Row 2
$ws.Cells.Item(2,4).Value = 44187
$ws.Cells.Item(2,11).Value = "Dina"
$ws.Cells.Item(2,13).Value = 120
$ws.Cells.Item(2,14).Value = 16000
$ws.Cells.Item(2,15).Value = 16000
$ws.Cells.Item(2,16).Value = 16000
$ws.Cells.Item(2,17).Value = "$/caja 18 kilos"
$ws.Cells.Item(2,19).Value = 889
$ws.Cells.Item(2,20).Value = 18

# Row 3
$ws.Cells.Item(3,4).Value = 44174
$ws.Cells.Item(3,12).Value = "Especial"
$ws.Cells.Item(3,13).Value = 200
$ws.Cells.Item(3,17).Value = "$/caja 10 kilos"
$ws.Cells.Item(3,18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(3,19).Value = 1500
$ws.Cells.Item(3,20).Value = 10

# Row 4
$ws.Cells.Item(4,4).Value = 44167
$ws.Cells.Item(4,11).Value = "Castle Brite"
$ws.Cells.Item(4,13).Value = 300
$ws.Cells.Item(4,17).Value = "$/caja 16 kilos granel"
$ws.Cells.Item(4,18).Value = "Provincia de Limarí"
$ws.Cells.Item(4,19).Value = 938
$ws.Cells.Item(4,20).Value = 16

# Row 5
$ws.Cells.Item(5,4).Value = 44176
$ws.Cells.Item(5,11).Value = "Castle Brite"
$ws.Cells.Item(5,13).Value = 100
$ws.Cells.Item(5,17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(5,18).Value = "Provincia de Limarí"

# Row 6
$ws.Cells.Item(6,4).Value = 44161
$ws.Cells.Item(6,13).Value = 150
$ws.Cells.Item(6,14).Value = 20000
$ws.Cells.Item(6,15).Value = 20000
$ws.Cells.Item(6,16).Value = 20000
$ws.Cells.Item(6,19).Value = 1111

# Row 7
$ws.Cells.Item(7,4).Value = 44162
$ws.Cells.Item(7,11).Value = "Castle Brite"
$ws.Cells.Item(7,13).Value = 200
$ws.Cells.Item(7,14).Value = 17000
$ws.Cells.Item(7,15).Value = 17000
$ws.Cells.Item(7,16).Value = 17000
$ws.Cells.Item(7,17).Value = "$/caja 16 kilos granel"
$ws.Cells.Item(7,18).Value = "Provincia de Limarí"
$ws.Cells.Item(7,19).Value = 1062
$ws.Cells.Item(7,20).Value = 16

# Row 8
$ws.Cells.Item(8,4).Value = 44162
$ws.Cells.Item(8,11).Value = "Castle Brite"
$ws.Cells.Item(8,12).Value = "Segunda"
$ws.Cells.Item(8,13).Value = 100
$ws.Cells.Item(8,14).Value = 15000
$ws.Cells.Item(8,15).Value = 15000
$ws.Cells.Item(8,16).Value = 15000
$ws.Cells.Item(8,17).Value = "$/caja 16 kilos granel"
$ws.Cells.Item(8,19).Value = 938
$ws.Cells.Item(8,20).Value = 16

# Row 9
$ws.Cells.Item(9,4).Value = 44160
$ws.Cells.Item(9,11).Value = "Dina"
$ws.Cells.Item(9,13).Value = 200
$ws.Cells.Item(9,14).Value = 20000
$ws.Cells.Item(9,15).Value = 20000
$ws.Cells.Item(9,16).Value = 20000
$ws.Cells.Item(9,17).Value = "$/caja 15 kilos"
$ws.Cells.Item(9,19).Value = 1333
$ws.Cells.Item(9,20).Value = 15

# Row 10
$ws.Cells.Item(10,4).Value = 44189
$ws.Cells.Item(10,11).Value = "Dina"
$ws.Cells.Item(10,12).Value = "Primera"
$ws.Cells.Item(10,13).Value = 50
$ws.Cells.Item(10,17).Value = "$/caja 18 kilos"
$ws.Cells.Item(10,19).Value = 833
$ws.Cells.Item(10,20).Value = 18

# Row 11
$ws.Cells.Item(11,4).Value = 44172
$ws.Cells.Item(11,12).Value = "Especial"
$ws.Cells.Item(11,13).Value = 120
$ws.Cells.Item(11,14).Value = 15000
$ws.Cells.Item(11,15).Value = 15000
$ws.Cells.Item(11,16).Value = 15000
$ws.Cells.Item(11,19).Value = 1500

# Row 12
$ws.Cells.Item(12,4).Value = 44172
$ws.Cells.Item(12,13).Value = 250
$ws.Cells.Item(12,14).Value = 11000
$ws.Cells.Item(12,15).Value = 11000
$ws.Cells.Item(12,16).Value = 11000
$ws.Cells.Item(12,17).Value = "$/caja 10 kilos"
$ws.Cells.Item(12,18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(12,19).Value = 1100
$ws.Cells.Item(12,20).Value = 10

# Row 13
$ws.Cells.Item(13,4).Value = 44186
$ws.Cells.Item(13,11).Value = "Dina"
$ws.Cells.Item(13,14).Value = 15000
$ws.Cells.Item(13,15).Value = 15000
$ws.Cells.Item(13,16).Value = 15000
$ws.Cells.Item(13,17).Value = "$/caja 18 kilos"
$ws.Cells.Item(13,18).Value = "Región Metropolitana"
$ws.Cells.Item(13,19).Value = 833

# Row 14
$ws.Cells.Item(14,4).Value = 44168
$ws.Cells.Item(14,11).Value = "Castle Brite"
$ws.Cells.Item(14,13).Value = 250
$ws.Cells.Item(14,14).Value = 10000
$ws.Cells.Item(14,15).Value = 10000
$ws.Cells.Item(14,16).Value = 10000
$ws.Cells.Item(14,17).Value = "$/caja 10 kilos"
$ws.Cells.Item(14,18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(14,19).Value = 1000
$ws.Cells.Item(14,20).Value = 10

# Row 15
$ws.Cells.Item(15,4).Value = 44168
$ws.Cells.Item(15,12).Value = "Primera"
$ws.Cells.Item(15,13).Value = 100
$ws.Cells.Item(15,14).Value = 17000
$ws.Cells.Item(15,15).Value = 17000
$ws.Cells.Item(15,16).Value = 17000
$ws.Cells.Item(15,17).Value = "$/caja 18 kilos"
$ws.Cells.Item(15,18).Value = "Provincia de Limarí"
$ws.Cells.Item(15,19).Value = 944
$ws.Cells.Item(15,20).Value = 18

# Row 16
$ws.Cells.Item(16,4).Value = 44181
$ws.Cells.Item(16,11).Value = "Dina"
$ws.Cells.Item(16,13).Value = 220
$ws.Cells.Item(16,14).Value = 17000
$ws.Cells.Item(16,15).Value = 17000
$ws.Cells.Item(16,16).Value = 17000
$ws.Cells.Item(16,17).Value = "$/caja 18 kilos"
$ws.Cells.Item(16,19).Value = 944
$ws.Cells.Item(16,20).Value = 18

# Row 17
$ws.Cells.Item(17,4).Value = 44179
$ws.Cells.Item(17,13).Value = 150
$ws.Cells.Item(17,14).Value = 18000
$ws.Cells.Item(17,15).Value = 18000
$ws.Cells.Item(17,16).Value = 18000
$ws.Cells.Item(17,18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(17,19).Value = 1000
